$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "32×38="
$t.Cell(1, 2).Range.Text = "18×85="
$t.Cell(1, 3).Range.Text = "69×67="
$t.Cell(1, 4).Range.Text = "18×31="
$t.Cell(1, 5).Range.Text = "13×94="
$t.Cell(5, 1).Range.Text = "71×39="
$t.Cell(5, 2).Range.Text = "34×23="
$t.Cell(5, 3).Range.Text = "83×71="
$t.Cell(5, 4).Range.Text = "26×81="
$t.Cell(5, 5).Range.Text = "89×86="
$t.Cell(10, 1).Range.Text = "46×88="
$t.Cell(10, 2).Range.Text = "58×31="
$t.Cell(10, 3).Range.Text = "15×18="
$t.Cell(10, 4).Range.Text = "78×89="
$t.Cell(10, 5).Range.Text = "63×47="
$t.Cell(15, 1).Range.Text = "53×91="
$t.Cell(15, 2).Range.Text = "56×85="
$t.Cell(15, 3).Range.Text = "13×50="
$t.Cell(15, 4).Range.Text = "42×86="
$t.Cell(15, 5).Range.Text = "59×52="
$t.Cell(20, 1).Range.Text = "15×75="
$t.Cell(20, 2).Range.Text = "35×12="
$t.Cell(20, 3).Range.Text = "14×56="
$t.Cell(20, 4).Range.Text = "42×52="
$t.Cell(20, 5).Range.Text = "22×88="
